$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab to match the player name.
$ws.Name = "Mohammed Siraj"

# Protect the numeric-looking text columns (runs, balls, fours, sixes, sr)
# from being auto-coerced into real numbers when assigned via .Value -
# the source data keeps them as plain text (numberStoredAsText).
$ws.Range("E2:I5").NumberFormat = "@"

# The "states" cell is blank for rows 3-5; force a text format so an
# empty-string assignment still leaves a (empty) cell behind instead of
# Excel clearing/removing it outright.
$ws.Range("D3:D5").NumberFormat = "@"

# Header row (a new "matchNo" column is inserted before the old column A,
# shifting every other header one column to the right).
$ws.Range("A1").Value = "matchNo"
$ws.Range("B1").Value = "teamName"
$ws.Range("C1").Value = "batterName"
$ws.Range("D1").Value = "states"
$ws.Range("E1").Value = "runs"
$ws.Range("F1").Value = "balls"
$ws.Range("G1").Value = "fours"
$ws.Range("H1").Value = "sixes"
$ws.Range("I1").Value = "sr"
$ws.Range("J1").Value = "opponentTeamName"
$ws.Range("K1").Value = "venue"
$ws.Range("L1").Value = "date"
$ws.Range("M1").Value = "result"

# Row 2
$ws.Range("A2").Value = "31st"
$ws.Range("B2").Value = "Royal Challengers Bangalore"
$ws.Range("C2").Value = "Mohammed Siraj"
$ws.Range("D2").Value = "c Varun b Russell"
$ws.Range("E2").Value = "8"
$ws.Range("F2").Value = "10"
$ws.Range("G2").Value = "1"
$ws.Range("H2").Value = "0"
$ws.Range("I2").Value = "80.00"
$ws.Range("J2").Value = "Kolkata Knight Riders"
$ws.Range("K2").Value = "Abu Dhabi"
$ws.Range("L2").Value = "September 20"
$ws.Range("M2").Value = "KKR won by 9 wickets (with 60 balls remaining)"

# Row 3
$ws.Range("A3").Value = "19th"
$ws.Range("B3").Value = "Royal Challengers Bangalore"
$ws.Range("C3").Value = "Mohammed Siraj"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "12"
$ws.Range("F3").Value = "14"
$ws.Range("G3").Value = "0"
$ws.Range("H3").Value = "1"
$ws.Range("I3").Value = "85.71"
$ws.Range("J3").Value = "Chennai Super Kings"
$ws.Range("K3").Value = "Wankhede"
$ws.Range("L3").Value = "April 25"
$ws.Range("M3").Value = "Super Kings won by 69 runs"

# Row 4
$ws.Range("A4").Value = "26th"
$ws.Range("B4").Value = "Royal Challengers Bangalore"
$ws.Range("C4").Value = "Mohammed Siraj"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "0"
$ws.Range("F4").Value = "1"
$ws.Range("G4").Value = "0"
$ws.Range("H4").Value = "0"
$ws.Range("I4").Value = "0.00"
$ws.Range("J4").Value = "Punjab Kings"
$ws.Range("K4").Value = "Ahmedabad"
$ws.Range("L4").Value = "April 30"
$ws.Range("M4").Value = "Punjab Kings won by 34 runs"

# Row 5
$ws.Range("A5").Value = "1st"
$ws.Range("B5").Value = "Royal Challengers Bangalore"
$ws.Range("C5").Value = "Mohammed Siraj"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "0"
$ws.Range("F5").Value = "1"
$ws.Range("G5").Value = "0"
$ws.Range("H5").Value = "0"
$ws.Range("I5").Value = "0.00"
$ws.Range("J5").Value = "Mumbai Indians"
$ws.Range("K5").Value = "Chennai"
$ws.Range("L5").Value = "April 09"
$ws.Range("M5").Value = "RCB won by 2 wickets"
